$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price) values are plain-looking numeric text in the source data
# (e.g. "580.74", "0.0000100") but must be stored as literal text, matching
# the original inlineStr cells. Forcing NumberFormat to text before the
# assignment keeps Excel from reinterpreting the string as a float (which would
# introduce binary rounding noise or drop significant trailing zeros), and
# ClearFormats() afterwards drops the "number stored as text" style tag Excel
# auto-applies so the cell keeps no explicit style, just like the source file.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.882.75'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.101.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +5.33%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.74'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.75'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.14%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.096.52'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.27%  '
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.43'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.65%  '
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.32'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.17%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.613.46'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +5.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.883.83'
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.102.18'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.19'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '480.66'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +7.91%  '
$ws.Range("E22").Value = '  +2.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.51'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.99'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.19'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.93%  '
$ws.Range("E26").Value = '  +4.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.07'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.97'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.86'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0000100'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.992'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.95'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.26%  '
$ws.Range("E39").Value = '  +7.41%  '
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.66'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0361'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.838.64'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +5.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '385.88'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.94'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.96'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.94%  '
$ws.Range("E51").Value = '  +2.77%  '
